$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "42.348.14"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +1.26%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.272.00"
$c.Style = "Normal"

$ws.Range("E4").Value = "  +0.03%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "305.99"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.61%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "97.51"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +5.13%  "

$ws.Range("E7").Value = "  -0.14%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.492"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.21%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "35.79"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +9.67%  "

$ws.Range("E11").Value = "  +0.00%  "

$ws.Range("E12").Value = "  -1.03%  "

$ws.Range("E13").Value = "  -0.58%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "2.624.36"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.02%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "14.38"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.44%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "2.279.54"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.16%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.795"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.01%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "42.246.97"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.22%  "

$ws.Range("E19").Value = "  -0.67%  "

$ws.Range("E20").Value = "  +0.16%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.96"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.22%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "67.60"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.71%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "240.55"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.43%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.60"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.55%  "

$ws.Range("E25").Value = "  +1.01%  "

$ws.Range("E26").Value = "  -0.24%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "23.82"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.99%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "37.58"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +6.18%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.53"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("E30").Value = "  +1.84%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "159.94"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.47%  "

$ws.Range("E32").Value = "  +0.16%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.08%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.16"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +4.53%  "

$ws.Range("E35").Value = "  -0.26%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "17.10"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("E37").Value = "  -0.36%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.34"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.72%  "

$ws.Range("E39").Value = "  +1.60%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.115"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.42%  "

$ws.Range("E41").Value = "  +3.99%  "

$ws.Range("E42").Value = "  +14.30%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.996.49"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.48%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0286"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.17%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "18.86"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -3.49%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.95"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.32%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "9.99"
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "53.27"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.80%  "

$ws.Range("E49").Value = "  +0.08%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "72.23"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.10%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "91.61"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.29%  "
